# Update gh-pages output: refresh dates (dot -> dash format) and ticket counts.
$wb = $excel.ActiveWorkbook

function Set-DateText {
    param($range, [string]$text)
    # Force text storage so Excel's autodetect doesn't turn "YYYY-MM-DD" into
    # a date serial; restore the original "General" formatting afterwards so
    # no stray number-format carries over onto the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# --- 展览 (Exhibition) sheet ---
$ws1 = $wb.Worksheets.Item("展览")
Set-DateText $ws1.Range("B2") "2024-03-09"
$ws1.Range("F2").Value = 650
Set-DateText $ws1.Range("B3") "2024-03-16"
Set-DateText $ws1.Range("B4") "2024-03-16"
$ws1.Range("F4").Value = 1490
Set-DateText $ws1.Range("B5") "2024-03-30"
$ws1.Range("F5").Value = 696

# --- 演出 (Performance) sheet ---
$ws2 = $wb.Worksheets.Item("演出")
Set-DateText $ws2.Range("B2") "2024-03-30"

# --- 全部类型 (All types) sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
Set-DateText $ws4.Range("B2") "2024-03-09"
$ws4.Range("F2").Value = 650
Set-DateText $ws4.Range("B3") "2024-03-16"
Set-DateText $ws4.Range("B4") "2024-03-16"
$ws4.Range("F4").Value = 1490
Set-DateText $ws4.Range("B5") "2024-03-30"
Set-DateText $ws4.Range("B6") "2024-03-30"
$ws4.Range("F6").Value = 696
